$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.543.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.38%  '
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.80'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3931'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3998'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.523'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.75%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.46'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08732'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.200'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +6.65%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001316'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.574'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.690.12'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07047'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.58'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.854'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.04'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.538.63'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.046'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +7.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.326'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.28'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.211'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.99'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.518'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.881.61'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.084'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08531'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.263'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +10.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.26'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.948'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2704'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.39'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02744'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +9.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09013'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.469'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7622'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7148'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.37'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.97%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.512'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.202'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.87'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.324'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +7.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07983'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.54%  '
